$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "96.624.95"
$ws.Range("E2").Value = "  -0.64%  "
$ws.Range("D3").Value = "3.681.58"
$ws.Range("E3").Value = "  +2.25%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'241.50"
$ws.Range("E5").Value = "  -0.71%  "
$ws.Range("D6").Value = "1.87"
$ws.Range("E6").Value = "  +10.50%  "
$ws.Range("D7").Value = "664.73"
$ws.Range("E7").Value = "  +1.26%  "
$ws.Range("D8").Value = "0.423"
$ws.Range("E8").Value = "  +0.84%  "
$ws.Range("D9").Value = "1.09"
$ws.Range("E9").Value = "  +1.95%  "
$ws.Range("D10").Value = "'1.00"
$ws.Range("E10").Value = "  +0.01%  "
$ws.Range("D11").Value = "3.680.65"
$ws.Range("E11").Value = "  +2.30%  "
$ws.Range("D12").Value = "45.97"
$ws.Range("E12").Value = "  +4.72%  "
$ws.Range("E13").Value = "  +0.43%  "
$ws.Range("D14").Value = "6.82"
$ws.Range("E14").Value = "  +5.70%  "
$ws.Range("D15").Value = "4.366.93"
$ws.Range("E15").Value = "  +2.24%  "
$ws.Range("D16").Value = "'0.0000269"
$ws.Range("E16").Value = "  +3.45%  "
$ws.Range("D17").Value = "96.429.34"
$ws.Range("E17").Value = "  -0.56%  "
$ws.Range("D18").Value = "8.95"
$ws.Range("E18").Value = "  +15.08%  "
$ws.Range("D19").Value = "3.683.59"
$ws.Range("E19").Value = "  +2.06%  "
$ws.Range("D20").Value = "12.92"
$ws.Range("E20").Value = "  +1.21%  "
$ws.Range("D21").Value = "18.69"
$ws.Range("E21").Value = "  +3.48%  "
$ws.Range("D22").Value = "0.528"
$ws.Range("E22").Value = "  -0.77%  "
$ws.Range("D23").Value = "527.63"
$ws.Range("E23").Value = "  +2.76%  "
$ws.Range("D24").Value = "3.46"
$ws.Range("E24").Value = "  +0.96%  "
$ws.Range("B25").Value = "PEPE"
$ws.Range("C25").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D25").Value = "'0.0000204"
$ws.Range("E25").Value = "  -0.30%  "
$ws.Range("B26").Value = "NEARProtocol"
$ws.Range("C26").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D26").Value = "7.06"
$ws.Range("E26").Value = "  +2.56%  "
$ws.Range("D27").Value = "102.33"
$ws.Range("E27").Value = "  +3.64%  "
$ws.Range("D28").Value = "13.17"
$ws.Range("E28").Value = "  +0.45%  "
$ws.Range("D29").Value = "3.877.52"
$ws.Range("E29").Value = "  +2.20%  "
$ws.Range("E30").Value = "  +8.95%  "
$ws.Range("D31").Value = "12.61"
$ws.Range("E31").Value = "  +7.21%  "
$ws.Range("D32").Value = "3.06"
$ws.Range("E32").Value = "  +0.95%  "
$ws.Range("D33").Value = "0.999"
$ws.Range("E33").Value = "  +0.07%  "
$ws.Range("E34").Value = "  +17.75%  "
$ws.Range("E35").Value = "  -0.20%  "
$ws.Range("D36").Value = "32.68"
$ws.Range("E36").Value = "  +2.98%  "
$ws.Range("D37").Value = "'1.00"
$ws.Range("E37").Value = "  +0.19%  "
$ws.Range("D38").Value = "660.15"
$ws.Range("E38").Value = "  +6.25%  "
$ws.Range("D39").Value = "0.593"
$ws.Range("E39").Value = "  +3.68%  "
$ws.Range("E40").Value = "  +0.73%  "
$ws.Range("D41").Value = "43.88"
$ws.Range("E41").Value = "  +32.12%  "
$ws.Range("E42").Value = "  +5.26%  "
$ws.Range("B43").Value = "ImmutableX"
$ws.Range("C43").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D43").Value = "'2.00"
$ws.Range("E43").Value = "  +3.63%  "
$ws.Range("B44").Value = "ARBITRUM"
$ws.Range("C44").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D44").Value = "0.969"
$ws.Range("E44").Value = "  +4.34%  "
$ws.Range("D45").Value = "6.52"
$ws.Range("E45").Value = "  +8.90%  "
$ws.Range("E46").Value = "  -0.02%  "
$ws.Range("D47").Value = "0.0467"
$ws.Range("E47").Value = "  +6.48%  "
$ws.Range("D48").Value = "0.449"
$ws.Range("E48").Value = "  +16.12%  "
$ws.Range("D49").Value = "2.32"
$ws.Range("E49").Value = "  +0.11%  "
$ws.Range("D50").Value = "3.68"
$ws.Range("E50").Value = "  +4.40%  "
$ws.Range("E51").Value = "  -0.16%  "
